# Update the "last_edited_time" timestamps for the rows that were touched
# by the report (rows 3, 4, 5, 7 and 13) and refresh the calculated
# "Thang 7" (row 5) metrics that come from the underlying Notion formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- last_edited_time (column D) updates ---------------------------------
$newTimestamp = "2024-07-17T12:15:00.000Z"

$ws.Range("D3").Value = $newTimestamp
$ws.Range("D4").Value = $newTimestamp
$ws.Range("D5").Value = $newTimestamp
$ws.Range("D7").Value = $newTimestamp
$ws.Range("D13").Value = $newTimestamp

# --- recalculated values for row 5 (Thang 7) ------------------------------
$ws.Range("T5").Value = 3500000
$ws.Range("W5").Value = 16537000
$ws.Range("AA5").Value = 21413000
$ws.Range("AE5").Value = 37950000
$ws.Range("AH5").Value = 34650000
$ws.Range("AK5").Value = 10
$ws.Range("AN5").Value = 3300000
$ws.Range("AQ5").Value = 38150000
